# Update cryptos list (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.121.79'
$ws.Cells.Item(2, 5).Value = '  +0.40%  '
$ws.Cells.Item(3, 4).Value = '1.917.98'
$ws.Cells.Item(3, 5).Value = '  +2.49%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '320.02'
$ws.Cells.Item(5, 5).Value = '  +0.16%  '
$ws.Cells.Item(6, 5).Value = '  +0.07%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5060'
$ws.Cells.Item(7, 5).Value = '  -0.41%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4073'
$ws.Cells.Item(8, 5).Value = '  +3.54%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08341'
$ws.Cells.Item(9, 5).Value = '  +1.78%  '
$ws.Cells.Item(10, 2).Value = 'OKB'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '42.35'
$ws.Cells.Item(10, 5).Value = '  +0.38%  '
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.109'
$ws.Cells.Item(11, 5).Value = '  +1.51%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '24.00'
$ws.Cells.Item(12, 5).Value = '  +5.17%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.916.60'
$ws.Cells.Item(13, 5).Value = '  +3.09%  '
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.421'
$ws.Cells.Item(14, 5).Value = '  +2.43%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.242'
$ws.Cells.Item(15, 5).Value = '  +1.05%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.004'
$ws.Cells.Item(16, 5).Value = '  +0.22%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '92.50'
$ws.Cells.Item(17, 5).Value = '  +0.65%  '
$ws.Cells.Item(18, 5).Value = '  +1.13%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06514'
$ws.Cells.Item(19, 5).Value = '  +1.30%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.53'
$ws.Cells.Item(20, 5).Value = '  +3.46%  '
$ws.Cells.Item(21, 5).Value = '  +0.02%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.946'
$ws.Cells.Item(22, 5).Value = '  +2.36%  '
$ws.Cells.Item(23, 4).Value = '30.131.66'
$ws.Cells.Item(23, 5).Value = '  +0.48%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.195'
$ws.Cells.Item(25, 5).Value = '  +2.45%  '
$ws.Cells.Item(26, 4).Value = '2.133.21'
$ws.Cells.Item(26, 5).Value = '  +2.67%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '21.89'
$ws.Cells.Item(27, 5).Value = '  +4.51%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '162.50'
$ws.Cells.Item(28, 5).Value = '  +0.79%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.266'
$ws.Cells.Item(29, 5).Value = '  +1.20%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '129.00'
$ws.Cells.Item(30, 5).Value = '  +1.51%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.134'
$ws.Cells.Item(31, 5).Value = '  +7.46%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.1045'
$ws.Cells.Item(32, 5).Value = '  +1.10%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.955'
$ws.Cells.Item(33, 5).Value = '  +0.84%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.787'
$ws.Cells.Item(34, 5).Value = '  +1.20%  '
$ws.Cells.Item(35, 5).Value = '  +1.65%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.330'
$ws.Cells.Item(36, 5).Value = '  +1.31%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06439'
$ws.Cells.Item(37, 5).Value = '  +1.68%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2154'
$ws.Cells.Item(38, 5).Value = '  +0.51%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.6515'
$ws.Cells.Item(39, 5).Value = '  +3.35%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.198'
$ws.Cells.Item(40, 5).Value = '  +2.21%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '8.604'
$ws.Cells.Item(41, 5).Value = '  +1.12%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '11.43'
$ws.Cells.Item(42, 5).Value = '  +1.82%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.208'
$ws.Cells.Item(43, 5).Value = '  +0.71%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '13.39'
$ws.Cells.Item(44, 5).Value = '  +3.97%  '
$ws.Cells.Item(45, 2).Value = 'NEARProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.189'
$ws.Cells.Item(45, 5).Value = '  +9.68%  '
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.6063'
$ws.Cells.Item(46, 5).Value = '  +2.69%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.626'
$ws.Cells.Item(47, 5).Value = '  -0.24%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.212'
$ws.Cells.Item(48, 5).Value = '  +0.78%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '122.30'
$ws.Cells.Item(49, 5).Value = '  -0.25%  '
$ws.Cells.Item(50, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.149'
$ws.Cells.Item(50, 5).Value = '  +1.67%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '79.00'
$ws.Cells.Item(51, 5).Value = '  +3.01%  '
